# "Change to new version" - replace the plaintext/MD5-hash password table
# with a new set of values, vertically-center the data cells, and extend
# the used range with new blank (but formatted) rows below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New Column A (plaintext) / Column D (MD5 hash) pairs for rows 2-11.
$newData = @(
    @(2,  "25OBEdIENt78",  "0a5840d9466c49b1a3a4fe3ff7ae6599"),
    @(3,  "aNNEx838",      "7f0d170908317ab2abe9cbb1c5448627"),
    @(4,  "CONCepT01",     "8e33ed543e50821e13ffa7454af27010"),
    @(5,  "8fIngEr",       "d804ecf5b670a04dde27c3670f0c1ac8"),
    @(6,  "121654River",   "cb6740a300d203fc1d7a9378825c2c2c"),
    @(7,  "sEntiENt13",    "4f108a47880629f180fbc3c0aa2d6fdf"),
    @(8,  "6poStER",       "9508d669d383d9526a31cf69dde6bde8"),
    @(9,  "8MiXTuRE5",     "869a174703dee3655aebd65402242a81"),
    @(10, "21385AgeNdA",   "aeb5c0e7f6f709c49cf81d6ba8791fd8"),
    @(11, "9805wHOlE",     "81ea785a5e780dc106921e3bcd6cbdfb")
)

foreach ($row in $newData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
}

# Vertically center the password/hash columns (A2:A11, D2:D11) - this is
# what turns the old "fill + border" style into the new "center" style.
$ws.Range("A2:A11").VerticalAlignment = -4108
$ws.Range("D2:D11").VerticalAlignment = -4108

# Extend formatting two rows below the table (row 12 stays blank/unused)
# down through row 22, matching the same vertical-center formatting.
$ws.Range("A13:D16").VerticalAlignment = -4108
$ws.Range("A17:A22").VerticalAlignment = -4108

# Move the active selection the way the author's session ended up.
$null = $ws.Range("C15").Select()

Write-Output "applied password table refresh"
